# Players Data.xlsx — "Logged Week 15 and simulated Week 16"
#
# Sheet "Rushing" (sheet1): a new player (R.Anderson) is logged into the
# roster ahead of the receivers who were previously listed starting at
# row 7, which pushes those names down one row each; a brand-new row 11
# is added at the bottom for I.Thomas. All of the 1DATT/2DATT/3DATT/RZATT
# counters are refreshed to the new week's cumulative totals.
#
# Sheet "Receiving" (sheet2): player list / order is unchanged; only the
# Short/Deep/RZ target & completion counters are refreshed.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Rushing sheet
# ---------------------------------------------------------------------
$rush = $wb.Worksheets.Item("Rushing")

# Row names shift down one slot (7->8->9->10) to make room for the newly
# logged R.Anderson at row 7; I.Thomas (previously row 10) becomes the
# brand-new row 11 at the bottom.
$rush.Range("B7").Value  = "R.Anderson"
$rush.Range("B8").Value  = "Dj.Moore"
$rush.Range("B9").Value  = "S.Smith"
$rush.Range("B10").Value = "T.Tremble"

# Refresh cumulative rushing totals: 1DATT, 2DATT, 3DATT, RZATT
$rush.Range("C2").Value = 21
$rush.Range("D2").Value = 21
$rush.Range("E2").Value = 17
$rush.Range("F2").Value = 9

$rush.Range("C4").Value = 31
$rush.Range("D4").Value = 16
$rush.Range("E4").Value = 5
$rush.Range("F4").Value = 8

$rush.Range("C6").Value = 14
$rush.Range("D6").Value = 11
$rush.Range("E6").Value = 4
$rush.Range("F6").Value = 3

$rush.Range("C7").Value = 2
$rush.Range("D7").Value = 0
$rush.Range("E7").Value = 0
$rush.Range("F7").Value = 1

$rush.Range("C8").Value = 1
$rush.Range("D8").Value = 2
$rush.Range("E8").Value = 1
$rush.Range("F8").Value = 0

$rush.Range("C9").Value = 0
$rush.Range("D9").Value = 1
$rush.Range("E9").Value = 0
$rush.Range("F9").Value = 0

$rush.Range("C10").Value = 0
$rush.Range("D10").Value = 1
$rush.Range("E10").Value = 0
$rush.Range("F10").Value = 1

# New row 11: I.Thomas — copy A10's formatting (bold/border/center style)
# onto A11 before filling in the row's values.
$rush.Range("A10").Copy($rush.Range("A11"))
$rush.Range("A11").Value = 9
$rush.Range("B11").Value = "I.Thomas"
$rush.Range("C11").Value = 1
$rush.Range("D11").Value = 0
$rush.Range("E11").Value = 0
$rush.Range("F11").Value = 0

# ---------------------------------------------------------------------
# Receiving sheet
# ---------------------------------------------------------------------
$rec = $wb.Worksheets.Item("Receiving")

$rec.Range("C2").Value = 16
$rec.Range("D2").Value = 10

$rec.Range("C4").Value = 26
$rec.Range("D4").Value = 18
$rec.Range("E4").Value = 2
$rec.Range("F4").Value = 1

$rec.Range("C5").Value = 70
$rec.Range("D5").Value = 35
$rec.Range("G5").Value = 6

$rec.Range("C6").Value = 98
$rec.Range("D6").Value = 64
$rec.Range("E6").Value = 35
$rec.Range("G6").Value = 12

$rec.Range("C7").Value = 21

$rec.Range("C8").Value = 17
$rec.Range("D8").Value = 14
$rec.Range("G8").Value = 3

$rec.Range("E11").Value = 1

$rec.Range("C13").Value = 22
$rec.Range("D13").Value = 15
$rec.Range("E13").Value = 6

$rec.Range("C14").Value = 23
$rec.Range("D14").Value = 12
